# Normalize the "Recorded By" (column G) entries: move every "System" /
# "system" token to the end of the comma-separated list, keeping the
# relative order of the remaining (non-System) entries, and keeping the
# relative order of the System-token(s) among themselves.
#
# e.g. "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#      "System, system, backup@backdoor.com"   -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $current = $cell.Value2

    if ($current -eq $null) { continue }

    $text = [string]$current
    if ($text -eq "") { continue }
    if ($text.IndexOf(",") -lt 0) { continue }

    $parts = $text -split ", "

    $others = @()
    $systems = @()
    foreach ($part in $parts) {
        if ($part.ToLower() -eq "system") {
            $systems += $part
        } else {
            $others += $part
        }
    }

    if ($systems.Count -eq 0 -or $others.Count -eq 0) { continue }

    $newParts = $others + $systems
    $newText = [string]::Join(", ", $newParts)

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
